$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.5523435399542
$ws.Range("C2").Value = 17.55742921144385
$ws.Range("D2").Value = 11.96228717994154

$ws.Range("B3").Value = 0.3649964739644579
$ws.Range("C3").Value = 0.2624334960621995
$ws.Range("D3").Value = 0.2870668474823946

$ws.Range("B4").Value = 0.3485797187544966
$ws.Range("C4").Value = 0.6427828021627909
$ws.Range("D4").Value = 0.2556116864232861

$ws.Range("B5").Value = 0.2851189159817535
$ws.Range("C5").Value = 0.4085519896680884
$ws.Range("D5").Value = 0.2528149115790507

$ws.Range("B6").Value = 4.921550602529484
$ws.Range("C6").Value = 4.136497483818174
$ws.Range("D6").Value = 4.574839854008596

$ws.Range("A7").Value = 0.5019609516481496
$ws.Range("B7").Value = 0.4312174909554408
$ws.Range("C7").Value = 0.4839749065006404
$ws.Range("D7").Value = 0.2327944710257013
